$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update crypto price/volume data cells to refreshed values.
# D-column (Price) values are forced as text via a leading apostrophe
# (matching the source inlineStr cell type) and the style is reset to
# "Normal" afterwards so no stray number-format style gets attached,
# since some numeric-looking price strings would otherwise be silently
# reinterpreted by Excel as actual numbers (losing formatting).

$ws.Range("D2").Value = "'65.962.63"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +1.83%  "

$ws.Range("D3").Value = "'3.498.05"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +0.58%  "

$ws.Range("E4").Value = "  -0.01%  "

$ws.Range("D5").Value = "'583.09"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.19%  "

$ws.Range("D6").Value = "'162.97"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +3.28%  "

$ws.Range("E7").Value = "  -0.01%  "

$ws.Range("D8").Value = "'3.503.28"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +0.57%  "

$ws.Range("D9").Value = "'0.586"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +6.43%  "

$ws.Range("D10").Value = "'7.34"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -3.32%  "

$ws.Range("D11").Value = "'0.127"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +1.12%  "

$ws.Range("D12").Value = "'0.448"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +0.07%  "

$ws.Range("D13").Value = "'4.106.02"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +0.75%  "

$ws.Range("D14").Value = "'0.135"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -1.43%  "

$ws.Range("D15").Value = "'0.0000198"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -0.07%  "

$ws.Range("D16").Value = "'28.79"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +3.64%  "

$ws.Range("D17").Value = "'65.950.77"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +1.73%  "

$ws.Range("D18").Value = "'3.468.53"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +0.26%  "

$ws.Range("D19").Value = "'6.48"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.04%  "

$ws.Range("D20").Value = "'14.44"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -0.02%  "

$ws.Range("D21").Value = "'396.73"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -0.38%  "

$ws.Range("D22").Value = "'8.36"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -2.99%  "

$ws.Range("D23").Value = "'74.16"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +1.71%  "

$ws.Range("D24").Value = "'0.552"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.56%  "

$ws.Range("E25").Value = "  +0.35%  "

$ws.Range("D26").Value = "'0.0000127"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +4.49%  "

$ws.Range("D27").Value = "'9.67"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.40%  "

$ws.Range("E28").Value = "  +0.40%  "

$ws.Range("E29").Value = "  +0.27%  "

$ws.Range("D30").Value = "'6.44"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +8.53%  "

$ws.Range("D31").Value = "'1.46"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +5.78%  "

$ws.Range("D32").Value = "'2.07"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +0.51%  "

$ws.Range("D33").Value = "'6.69"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +0.82%  "

$ws.Range("D34").Value = "'23.84"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -0.16%  "

$ws.Range("E35").Value = "  -0.02%  "

$ws.Range("D36").Value = "'7.19"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +3.01%  "

$ws.Range("D37").Value = "'1.54"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +2.25%  "

$ws.Range("D38").Value = "'162.35"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +1.24%  "

$ws.Range("D39").Value = "'2.00"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +5.04%  "

$ws.Range("D40").Value = "'3.030.21"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +4.07%  "

$ws.Range("D41").Value = "'0.0778"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -0.56%  "

$ws.Range("D42").Value = "'27.43"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -2.07%  "

$ws.Range("B43").Value = "VeChain"
$ws.Range("C43").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D43").Value = "'0.0325"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +0.70%  "

$ws.Range("B44").Value = "Filecoin"
$ws.Range("C44").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D44").Value = "'4.57"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +2.93%  "

$ws.Range("D45").Value = "'42.90"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +3.01%  "

$ws.Range("D46").Value = "'0.781"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +0.09%  "

$ws.Range("D47").Value = "'25.98"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +12.60%  "

$ws.Range("D48").Value = "'1.13"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +1.86%  "

$ws.Range("D49").Value = "'2.29"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +5.17%  "

$ws.Range("B50").Value = "Bittensor"
$ws.Range("C50").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D50").Value = "'316.35"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +5.96%  "

$ws.Range("D51").Value = "'6.72"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +2.83%  "
